$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N7").Value = "done"
$ws.Range("N8").Value = "done"
$ws.Range("N9").Value = "done"
